# Refresh cryptocurrency prices and 1h volume-change percentages
# (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.496.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '''1.729.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('D4').Value = '''0.9987'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''246.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').Value = '''0.9994'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.4798'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.69%  '
$ws.Range('D8').Value = '''0.2688'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').Value = '''0.06227'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '''1.728.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').Value = '''0.07136'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('D13').Value = '''0.6196'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.78%  '
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('D15').Value = '''77.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '''0.9993'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '''26.503.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '''0.9991'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').Value = '''0.000006963'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '''1.950.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').Value = '''4.541'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = '''8.966'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').Value = '''5.299'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('D25').Value = '''136.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').Value = '''15.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('D27').Value = '''1.805'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.12%  '
$ws.Range('D28').Value = '''1.407'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').Value = '''106.89'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = '''3.987'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').Value = '''0.08033'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').Value = '''0.04573'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('D34').Value = '''0.9987'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = '''2.617'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('D36').Value = '''0.6381'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.32%  '
$ws.Range('D37').Value = '''0.9910'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('D38').Value = '''0.9354'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.01%  '
$ws.Range('D39').Value = '''2.085'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.59%  '
$ws.Range('D40').Value = '''2.414'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').Value = '''105.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.62%  '
$ws.Range('D42').Value = '''1.007'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '''5.716'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.04%  '
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('D45').Value = '''0.3914'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.63%  '
$ws.Range('D46').Value = '''6.980'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.60%  '
$ws.Range('E47').Value = '  +4.04%  '
$ws.Range('D48').Value = '''0.05324'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '''31.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = '''7.946'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.75%  '
$ws.Range('E51').Value = '  +3.65%  '
